$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 1084
$ws1.Range("G3").Value = 65
$ws1.Range("G4").Value = "不可售"
$ws1.Range("F5").Value = 3342
$ws1.Range("F7").Value = 301
$ws1.Range("F8").Value = 37
$ws1.Range("F10").Value = 18
$ws1.Range("F12").Value = 111
$ws1.Range("F13").Value = 197
$ws1.Range("F14").Value = 34
$ws1.Range("F15").Value = 86
$ws1.Range("F16").Value = 2809
$ws1.Range("F17").Value = 1116

# Update "全部类型" sheet (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F4").Value = 1084
$ws4.Range("G4").Value = 65
$ws4.Range("G5").Value = "不可售"
$ws4.Range("F6").Value = 3342
$ws4.Range("F8").Value = 301
$ws4.Range("F9").Value = 37
$ws4.Range("F12").Value = 18
$ws4.Range("F14").Value = 111
$ws4.Range("F15").Value = 197
$ws4.Range("F16").Value = 34
$ws4.Range("F17").Value = 86
$ws4.Range("F18").Value = 2809
$ws4.Range("F19").Value = 1116
